$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (set Y1 "label" before X1 "classifier 1classifier 2" so the
# shared-string table is built in the same order as the target file:
# index 0 = "label", index 1 = "classifier 1classifier 2").
$ws.Range("Y1").Value = "label"
$ws.Range("X1").Value = "classifier 1classifier 2"

# Data rows (columns W, X, Y = 23, 24, 25)
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 0.41
$ws.Range("Y2").Value = 1

$ws.Range("W3").Value = 0.72
$ws.Range("X3").Value = 0
$ws.Range("Y3").Value = 0

$ws.Range("W4").Value = 0.99
$ws.Range("X4").Value = 0
$ws.Range("Y4").Value = 0

$ws.Range("W5").Value = 0.14
$ws.Range("X5").Value = 0.73
$ws.Range("Y5").Value = 1

$ws.Range("W6").Value = 0
$ws.Range("X6").Value = 0.62
$ws.Range("Y6").Value = 0

$ws.Range("W7").Value = 0.94
$ws.Range("X7").Value = 1
$ws.Range("Y7").Value = 1

$ws.Range("W8").Value = 0.1
$ws.Range("X8").Value = 1
$ws.Range("Y8").Value = 0

$ws.Range("W9").Value = 0.77
$ws.Range("X9").Value = 0.14
$ws.Range("Y9").Value = 1

$ws.Range("W10").Value = 0.02
$ws.Range("X10").Value = 0
$ws.Range("Y10").Value = 0

$ws.Range("W11").Value = 1
$ws.Range("X11").Value = 0.55
$ws.Range("Y11").Value = 1

# Match the saved view state: T13 is the active/selected cell.
$ws.Range("T13").Select() | Out-Null
